$wb = $excel.ActiveWorkbook

# This edit regenerates the localization-status report for a new source
# file (cf7b8d21-28c7-4824-8621-cd8799ab6fde.md, replacing
# 743b15eb-139e-47e3-840b-579365822d91.md) and marks it ready for a fresh
# handoff: new handoff xliff names/timestamps, and the handback
# file/timestamp cleared back to "not yet received" (0001-01-01 00:00:00).

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "cf7b8d21-28c7-4824-8621-cd8799ab6fde.md"
$ws1.Range("B2").Value = "e2e\cf7b8d21-28c7-4824-8621-cd8799ab6fde.md"
$ws1.Range("G2").Value = "2016-08-18 11:03:08"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5939d601455e826eae4063ee8c0afe16e61ba3d6/e2e/743b15eb-139e-47e3-840b-579365822d91.md", "", "", "e2e\cf7b8d21-28c7-4824-8621-cd8799ab6fde.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "cf7b8d21-28c7-4824-8621-cd8799ab6fde.md"
$ws2.Range("G2").Value = "cf7b8d21-28c7-4824-8621-cd8799ab6fde.100a6635f79bb3f3b53f5bd84c8f4e4ca77e0c4e.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-18 11:02:58"
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

# Latest Target File / Latest Handback File no longer point at a handed
# back file, so the values (and their hyperlink) are cleared.
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5939d601455e826eae4063ee8c0afe16e61ba3d6/e2e/743b15eb-139e-47e3-840b-579365822d91.md", "", "", "cf7b8d21-28c7-4824-8621-cd8799ab6fde.md")
$ws2.Range("I2").Font.Underline = $false

$ws2.Columns.Item(9).ColumnWidth = 17.75
$ws2.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "cf7b8d21-28c7-4824-8621-cd8799ab6fde.md"
$ws3.Range("G2").Value = "cf7b8d21-28c7-4824-8621-cd8799ab6fde.100a6635f79bb3f3b53f5bd84c8f4e4ca77e0c4e.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-18 11:03:08"
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5939d601455e826eae4063ee8c0afe16e61ba3d6/e2e/743b15eb-139e-47e3-840b-579365822d91.md", "", "", "cf7b8d21-28c7-4824-8621-cd8799ab6fde.md")
$ws3.Range("I2").Font.Underline = $false

$ws3.Columns.Item(9).ColumnWidth = 17.75
$ws3.Columns.Item(10).ColumnWidth = 20.8
